$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E - shifts old E:H to F:I
$ws.Columns("E").Insert()

# Populate the new column E with the new test step (submit / id=password)
$ws.Range("E1").Value = "submit"
$ws.Range("E2").Value = "id=password"

# E3 stays empty but needs the same cell style as the JSON-snippet cells (H3)
$ws.Range("H3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match column E's width to its neighbouring "id=*" columns (~13.125 chars)
$ws.Columns("E").ColumnWidth = 12.43

# Move the active selection to F7 (matches the saved view state)
$ws.Range("F7").Select()
